$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 61 - Hash Set class
$ws.Range("C61").Value = "List of arrays, to which values are index-assigned based on hashing of input integer"
$ws.Range("D61").Value = "O(1)"
$ws.Range("E61").Value = "Linear"
$ws.Range("F61").Value = 45888
$ws.Range("G61").Value = "12 minutes"
$ws.Range("H61").Value = "20 minutes"
$ws.Range("I61").Value = "Got it mostly right, but forgot to consider resizing (for better retrieval)"
$ws.Range("J61").Value = "Almost"
$ws.Range("K61").Value = "That it was pretty similar to the dynamic array problem"
$ws.Range("L61").Value = "[] * len gave me copies of the same object, use list comprehension instead"
$ws.Range("M61").Value = "No"
$ws.Range("N61").Value = "Revise this implementation"
$ws.Range("O61").Value = "List comprehension to declare array of empty arrays"
$ws.Range("P61").Value = 3
$ws.Range("Q61").Value = 3
$ws.Range("R61").Value = 3
$ws.Range("S61").Value = 3
$ws.Rows.Item(61).RowHeight = 60

# Row 62 - Hash Set class extensions
$ws.Range("C62").Value = "Iterate through buckets to generate list, and use that as basis for union and extend"
$ws.Range("D62").Value = "Linear"
$ws.Range("E62").Value = "Linear"
$ws.Range("F62").Value = 45888
$ws.Range("G62").Value = "12 minutes"
$ws.Range("H62").Value = "20 minutes"
$ws.Range("I62").Value = "Got this one almost perfectly right"
$ws.Range("J62").Value = "Yes"
$ws.Range("K62").Value = "That I needed to iterate, and that involved traversing the buckets"
$ws.Range("L62").Value = "Appending instead of extending, extending is faster"
$ws.Range("M62").Value = "No"
$ws.Range("N62").Value = "Good to know"
$ws.Range("O62").Value = "Load factor logic"
$ws.Range("P62").Value = 4
$ws.Range("Q62").Value = 4
$ws.Range("R62").Value = 4
$ws.Range("S62").Value = 4
$ws.Rows.Item(62).RowHeight = 60

# Row 63 - Multiset
$ws.Range("C63").Value = "Minor changes to has set implementation"
$ws.Range("D63").Value = "O(1)"
$ws.Range("E63").Value = "Linear"
$ws.Range("F63").Value = 45888
$ws.Range("G63").Value = "12 minutes"
$ws.Range("H63").Value = "20 minutes"
$ws.Range("I63:N63").Merge()
$ws.Range("I63").Value = "All good"
$ws.Range("O63").Value = "Knuth multiplicative and base128 algorithms"
$ws.Range("P63").Value = 4
$ws.Range("Q63").Value = 4
$ws.Range("R63").Value = 4
$ws.Range("S63").Value = 4
$ws.Rows.Item(63).RowHeight = 40
